# Update the "Report for Handback" timestamps that are regenerated each run.
$wb = $excel.ActiveWorkbook

# Sheet "Overview": Latest HO Xliff Generate Date for the first file (row 2, col G)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-24 09:09:14"

# Sheet "zh-cn": Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-24 09:09:03"
$wsZhCn.Range("K2").Value = "2016-08-24 09:09:31"

# Sheet "de-de": Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-24 09:09:14"
$wsDeDe.Range("K2").Value = "2016-08-24 09:09:38"
